# Apply the Alvearie FHIR IG deploy update (StructureDefinition-communication-contact)
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Metadata")
$ws2 = $wb.Worksheets.Item("Elements")

# --- Metadata sheet ---
# Version bump
$ws1.Range("B3").Value = "6.0.0"

# Date bump
$ws1.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value now populated
$ws1.Range("B9").Value = "Alvearie Team"

# Row 10 "Contact / No display for ContactDetail" becomes "Jurisdiction / United States of America"
$ws1.Range("A10").Value = "Jurisdiction"
$ws1.Range("B10").Value = "United States of America"

# Row 11 (duplicate "Contact / No display for ContactDetail") is removed entirely,
# rows below shift up by one
$ws1.Rows.Item(11).Delete()

# --- Elements sheet ---
# Root Extension row: Short / Definition now reflect the profile title/description
# instead of the generic placeholders
$ws2.Range("K2").Value = "Communication Contact"
$ws2.Range("L2").Value = "Phone number, email, or address to contact"
